$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New prediction rows for the "2021-01-09" block (Canada), extending the
# table from row 53 down to row 59 with one additional week
# ("14 Feb -- 20 Feb 2021") and its prediction value.

# Data for new rows 54..59:
#   Column A: the day the prediction was made -> same as row 48..53 ("2021-01-09")
#   Column B: the week of the target variable (label)
#   Column D: Prediction value
#   Column F: Model -> "KNN"
$weeks = @(
    "10 Jan -- 16 Jan 2021",
    "17 Jan -- 23 Jan 2021",
    "24 Jan -- 30 Jan 2021",
    "31 Jan -- 06 Feb 2021",
    "07 Feb -- 13 Feb 2021",
    "14 Feb -- 20 Feb 2021"
)
$preds = @(99.62, 106.65, 109.38, 105.69, 80.8, 101.37)

$startRow = 54
for ($i = 0; $i -lt $weeks.Length; $i++) {
    $r = $startRow + $i

    # Column A holds a date-like text ("2021-01-09"). Assigning a plain
    # string here would make Excel auto-convert it into a real date
    # serial number, so instead copy the already-typed text cell from the
    # existing block (row 53) which carries the exact same value and
    # cell formatting/type.
    $ws.Cells.Item(53, 1).Copy($ws.Cells.Item($r, 1))

    $ws.Cells.Item($r, 2).Value = $weeks[$i]
    $ws.Cells.Item($r, 4).Value = $preds[$i]
    $ws.Cells.Item($r, 6).Value = "KNN"
}
